$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1.53
$ws.Range("I3").Value = 5.5
$ws.Range("Z3").Value = 11
$ws.Range("AD3").Value = 9
$ws.Range("AW3").Value = 7.5
$ws.Range("AX3").Value = 29
